$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation needs to be inserted right before the
# current row 431. Inserting the row shifts rows 431-489 down to 432-490
# (carrying all of their existing values/formatting with them), leaving
# row 431 empty and ready to receive the new record.
$ws.Rows.Item(431).Insert()

# Fill in the new record in row 431.
$ws.Range("A431").Value = 8
$ws.Range("B431").Value = "Terminal La Palmera de La Serena"
$ws.Range("C431").Value = "Coquimbo"
$ws.Range("D431").Value = 44984
$ws.Range("E431").Value = 4
$ws.Range("F431").Value = 100114013
$ws.Range("G431").Value = "Zanahoria"
$ws.Range("H431").Value = "Sin especificar"
$ws.Range("I431").Value = "Primera"
$ws.Range("J431").Value = 440
$ws.Range("K431").Value = 5500
$ws.Range("L431").Value = 6000
$ws.Range("M431").Value = 5750
$ws.Range("N431").Value = "$/saco 20 kilos"
$ws.Range("O431").Value = "Provincia del Elquí"
$ws.Range("P431").Value = 288
$ws.Range("Q431").Value = 20
$ws.Range("R431").Value = "Hortaliza"
